$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$t = $m.Theme
$tcs = $t.ThemeColorScheme
$c1 = $tcs.Item(1)
$c1.RGB = 255
Write-Host "set done"
